$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value="57.863.38"},
    @{Cell="E2"; Value="  -1.69%  "},
    @{Cell="D3"; Value="2.453.20"},
    @{Cell="E3"; Value="  -1.72%  "},
    @{Cell="D4"; Value="1.00"},
    @{Cell="E4"; Value="  -0.26%  "},
    @{Cell="D5"; Value="517.20"},
    @{Cell="E5"; Value="  -3.59%  "},
    @{Cell="D6"; Value="131.86"},
    @{Cell="E6"; Value="  -3.01%  "},
    @{Cell="D7"; Value="0.999"},
    @{Cell="E7"; Value="  +0.02%  "},
    @{Cell="E8"; Value="  -1.82%  "},
    @{Cell="D9"; Value="2.457.65"},
    @{Cell="E9"; Value="  -2.64%  "},
    @{Cell="D10"; Value="0.0980"},
    @{Cell="E10"; Value="  -3.05%  "},
    @{Cell="E11"; Value="  -0.11%  "},
    @{Cell="E12"; Value="  -1.19%  "},
    @{Cell="E13"; Value="  -2.39%  "},
    @{Cell="D14"; Value="2.889.88"},
    @{Cell="E14"; Value="  -2.57%  "},
    @{Cell="D15"; Value="57.786.38"},
    @{Cell="E15"; Value="  -1.72%  "},
    @{Cell="D16"; Value="22.19"},
    @{Cell="E16"; Value="  -3.44%  "},
    @{Cell="E17"; Value="  -2.66%  "},
    @{Cell="D18"; Value="2.451.51"},
    @{Cell="E18"; Value="  -2.61%  "},
    @{Cell="D19"; Value="10.64"},
    @{Cell="E19"; Value="  -3.76%  "},
    @{Cell="B20"; Value="BitcoinCash"},
    @{Cell="C20"; Value="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"},
    @{Cell="D20"; Value="319.09"},
    @{Cell="E20"; Value="  -1.13%  "},
    @{Cell="B21"; Value="Polkadot"},
    @{Cell="C21"; Value="https://coinranking.com/coin/25W7FG7om+polkadot-dot"},
    @{Cell="D21"; Value="4.15"},
    @{Cell="E21"; Value="  -2.57%  "},
    @{Cell="E22"; Value="  +0.13%  "},
    @{Cell="D23"; Value="5.72"},
    @{Cell="E23"; Value="  -3.77%  "},
    @{Cell="D24"; Value="64.29"},
    @{Cell="E24"; Value="  -1.12%  "},
    @{Cell="E25"; Value="  -2.73%  "},
    @{Cell="D26"; Value="0.998"},
    @{Cell="E26"; Value="  +0.05%  "},
    @{Cell="E27"; Value="  -2.73%  "},
    @{Cell="D28"; Value="7.31"},
    @{Cell="E28"; Value="  -2.76%  "},
    @{Cell="E29"; Value="  -4.06%  "},
    @{Cell="D30"; Value="165.69"},
    @{Cell="E30"; Value="  -3.01%  "},
    @{Cell="E31"; Value="  -4.03%  "},
    @{Cell="D32"; Value="6.20"},
    @{Cell="E32"; Value="  -6.55%  "},
    @{Cell="D33"; Value="1.16"},
    @{Cell="E33"; Value="  -0.65%  "},
    @{Cell="E34"; Value="  +0.05%  "},
    @{Cell="E35"; Value="  +0.22%  "},
    @{Cell="E36"; Value="  -1.78%  "},
    @{Cell="E37"; Value="  -6.86%  "},
    @{Cell="E38"; Value="  -2.78%  "},
    @{Cell="E39"; Value="  -4.19%  "},
    @{Cell="B40"; Value="OKB"},
    @{Cell="C40"; Value="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"},
    @{Cell="D40"; Value="36.16"},
    @{Cell="E40"; Value="  -1.96%  "},
    @{Cell="B41"; Value="SuiNetwork"},
    @{Cell="C41"; Value="https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"},
    @{Cell="D41"; Value="0.785"},
    @{Cell="E41"; Value="  -3.12%  "},
    @{Cell="B42"; Value="Filecoin"},
    @{Cell="C42"; Value="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"},
    @{Cell="D42"; Value="3.42"},
    @{Cell="E42"; Value="  -4.53%  "},
    @{Cell="B43"; Value="Bittensor"},
    @{Cell="C43"; Value="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"},
    @{Cell="D43"; Value="270.71"},
    @{Cell="E43"; Value="  -4.82%  "},
    @{Cell="B44"; Value="RenderToken"},
    @{Cell="C44"; Value="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"},
    @{Cell="D44"; Value="5.00"},
    @{Cell="E44"; Value="  -3.03%  "},
    @{Cell="B45"; Value="Mantle"},
    @{Cell="C45"; Value="https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"},
    @{Cell="D45"; Value="0.588"},
    @{Cell="E45"; Value="  -2.99%  "},
    @{Cell="B46"; Value="Aave"},
    @{Cell="C46"; Value="https://coinranking.com/coin/ixgUfzmLR+aave-aave"},
    @{Cell="D46"; Value="124.47"},
    @{Cell="E46"; Value="  -4.07%  "},
    @{Cell="B47"; Value="Stellar"},
    @{Cell="C47"; Value="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"},
    @{Cell="D47"; Value="0.0905"},
    @{Cell="E47"; Value="  -1.93%  "},
    @{Cell="B48"; Value="Hedera"},
    @{Cell="C48"; Value="https://coinranking.com/coin/jad286TjB+hedera-hbar"},
    @{Cell="D48"; Value="0.0485"},
    @{Cell="E48"; Value="  -3.81%  "},
    @{Cell="B49"; Value="VeChain"},
    @{Cell="C49"; Value="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"},
    @{Cell="D49"; Value="0.0210"},
    @{Cell="E49"; Value="  -4.26%  "},
    @{Cell="B50"; Value="InjectiveProtocol"},
    @{Cell="C50"; Value="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"},
    @{Cell="D50"; Value="16.66"},
    @{Cell="E50"; Value="  -4.01%  "},
    @{Cell="B51"; Value="Maker"},
    @{Cell="C51"; Value="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"},
    @{Cell="D51"; Value="1.722.69"},
    @{Cell="E51"; Value="  -1.73%  "}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value2 = $u.Value
    $rng.Style = "Normal"
}
